$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E, shifting the old E column (DDL) to F
$ws.Range("E1").EntireColumn.Insert()

# New header cell E1: "ObjectType" — match the bold/bordered/centered header
# formatting already used by the other header cells (A1:D1, now also F1).
$ws.Range("E1").Value = "ObjectType"
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("E1").VerticalAlignment = -4160    # xlTop
$ws.Range("E1").Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$ws.Range("E1").Borders.Item(8).LineStyle = 1   # xlEdgeTop
$ws.Range("E1").Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$ws.Range("E1").Borders.Item(10).LineStyle = 1  # xlEdgeRight

# New data cell E2: "Sconosciuto" (plain, unstyled — matches the rest of row 2)
$ws.Range("E2").Value = "Sconosciuto"
